$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value = 10
$ws.Cells.Item(40, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(40, 3).Value = "La Araucanía"
$ws.Cells.Item(40, 4).Value = 44719
$ws.Cells.Item(40, 5).Value = 9
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100108
$ws.Cells.Item(40, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(40, 9).Value = 100108007
$ws.Cells.Item(40, 10).Value = "Coco"
$ws.Cells.Item(40, 11).Value = "Sin especificar"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 25
$ws.Cells.Item(40, 14).Value = 28000
$ws.Cells.Item(40, 15).Value = 28000
$ws.Cells.Item(40, 16).Value = 28000
$ws.Cells.Item(40, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(40, 18).Value = "Perú"
$ws.Cells.Item(40, 19).Value = 1400
$ws.Cells.Item(40, 20).Value = 20
